# GEORGIA_2016.xlsx style edit:
#  1) Rename header row (A1:D1) to short machine-friendly column names.
#  2) Title-case the Spanish connector words ("de", "del", "la", "las",
#     "los", "el", "y") whenever they occur as a non-leading word inside
#     a state/municipality name in columns A and B.
#  3) Fix a couple of 1-ULP floating point rounding artifacts in column D.
#  4) Drop the trailing footnote rows (1591-1596) so the sheet's used
#     range / dimension shrinks back to A1:D1590.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header row rename
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# ---------------------------------------------------------------------
# 2) Title-case connector words in columns A and B, rows 2..1590
# ---------------------------------------------------------------------
$connectors = @{ "de" = 1; "del" = 1; "la" = 1; "las" = 1; "los" = 1; "el" = 1; "y" = 1 }

function Fix-ConnectorCase($text) {
    $parts = $text.Split(" ")
    $changed = $false
    for ($i = 1; $i -lt $parts.Length; $i++) {
        $w = $parts[$i]
        if ($connectors.ContainsKey($w)) {
            $parts[$i] = $w.Substring(0, 1).ToUpper() + $w.Substring(1)
            $changed = $true
        }
    }
    if ($changed) {
        return [string]::Join(" ", $parts)
    }
    return $null
}

$lastRow = 1590
for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in 1, 2) {
        $cell = $ws.Cells.Item($r, $col)
        $v = $cell.Value()
        if ($v -ne $null -and $v.GetType().Name -eq "String" -and $v.Length -gt 0) {
            $fixed = Fix-ConnectorCase $v
            if ($fixed -ne $null) {
                $cell.Value = $fixed
            }
        }
    }
}

# ---------------------------------------------------------------------
# 3) Floating point 1-ULP fixes in column D
# ---------------------------------------------------------------------
$floatFixRows649 = 4, 102, 106, 235, 553, 673, 691, 702, 974, 1047, 1232, 1420, 1437, 1554
foreach ($r in $floatFixRows649) {
    $ws.Cells.Item($r, 4).Value = 0.0009624639076034648
}
$ws.Cells.Item(169, 4).Value = 0.009740134744947063

# ---------------------------------------------------------------------
# 4) Drop footnote rows 1591-1596 -> shrink used range to A1:D1590
# ---------------------------------------------------------------------
$ws.Range("A1591:D1596").Clear()
